$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update statistics for Reino Unido (row 9) ---
$ws.Range("B9").Value = 108692
$ws.Range("C9").Value = 5599
$ws.Range("E9").Value = 93772
$ws.Range("G9").Value = 847
$ws.Range("H9").Value = 14576

# --- Update statistics for Suiza (row 18) ---
$ws.Range("E18").Value = 9863
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = 1315

# --- Update statistics for Austria (row 20) ---
$ws.Range("B20").Value = 14553
$ws.Range("C20").Value = 77
$ws.Range("E20").Value = 4439

# --- Update statistics for Arabia Saudita (row 32) ---
$ws.Range("F32").Value = 74

# --- Update statistics for Finlandia (row 50) ---
$ws.Range("E50").Value = 1707
$ws.Range("F50").Value = 73
$ws.Range("G50").Value = 7
$ws.Range("H50").Value = 82

# --- Update statistics for Republica de Macedonia (row 77) ---
$ws.Range("B77").Value = 1117
$ws.Range("C77").Value = 36
$ws.Range("D77").Value = 139
$ws.Range("E77").Value = 929
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 49

# --- Swap Madagascar / Trinidad yTobago rows and update Madagascar's figures ---
# Row 134 currently holds Trinidad yTobago, row 135 currently holds Madagascar.
# After the edit, row 134 holds Madagascar (with updated figures) and row 135
# holds Trinidad yTobago (figures unchanged, just relocated).
$ws.Range("A134").Value = "Madagascar"
$ws.Range("B134").Value = 117
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 33
$ws.Range("E134").Value = 84
$ws.Range("F134").Value = 1
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0

$ws.Range("A135").Value = "Trinidad yTobago"
$ws.Range("B135").Value = 114
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 20
$ws.Range("E135").Value = 86
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 8
